# patients.xlsx — "add patient now working, delete is work in progress,
# tweaked database interfacing"
#
# The sheet used to hold 2 demo rows with a wide block of placeholder
# numeric columns (C:L), three duplicated date columns (M:O) and an "age
# group" column (Q). The new database-backed version only keeps:
#   A  -> row id
#   B  -> patient name (string)
#   M  -> created/updated timestamp (custom date-time format)
#   P  -> gender (string)
#   Q  -> age
# and now has 4 patient rows instead of 2 (a patient got added; the
# in-progress "delete" feature is not wired to this sheet yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Wipe the columns that the new database view no longer uses.
#    Reset to the "Normal" style first so no stray per-cell formatting
#    (and no orphan number formats) survive the clear.
# ------------------------------------------------------------------
$ws.Range("C1:L2").Style = "Normal"
$ws.Range("C1:L2").Value = $null

$ws.Range("N1:O2").Style = "Normal"
$ws.Range("N1:O2").Value = $null

$ws.Range("Q1:Q4").Style = "Normal"

# ------------------------------------------------------------------
# 2. Column M now carries a single "last touched" timestamp per row,
#    formatted as day/month/year hour:minute instead of the old
#    built-in datetime format.
# ------------------------------------------------------------------
$ws.Range("M1:M4").NumberFormat = "d/m/yy hh:mm"

# ------------------------------------------------------------------
# 3. Row 1 — existing patient, renamed "heba", age bumped to 44.
# ------------------------------------------------------------------
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "heba"
$ws.Range("M1").Value = 42457.05393122685
$ws.Range("P1").Value = "male"
$ws.Range("Q1").Value = 44

# ------------------------------------------------------------------
# 4. Row 2 — existing patient, renamed "ss", age reset to 1.
# ------------------------------------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "ss"
$ws.Range("M2").Value = 42457.05393173611
$ws.Range("P2").Value = "male"
$ws.Range("Q2").Value = 1

# ------------------------------------------------------------------
# 5. Row 3 — brand-new patient ("add patient" now working).
# ------------------------------------------------------------------
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "ss"
$ws.Range("M3").NumberFormat = "d/m/yy hh:mm"
$ws.Range("M3").Value = 42457.05386638889
$ws.Range("P3").Value = "male"
$ws.Range("Q3").Value = 1

# ------------------------------------------------------------------
# 6. Row 4 — another new patient.
# ------------------------------------------------------------------
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "ss"
$ws.Range("M4").NumberFormat = "d/m/yy hh:mm"
$ws.Range("M4").Value = 42457.04991763889
$ws.Range("P4").Value = "male"
$ws.Range("Q4").Value = 1

# ------------------------------------------------------------------
# 7. Column M is the only custom-width column left now.
# ------------------------------------------------------------------
$ws.Columns("M").ColumnWidth = 11.5

# ------------------------------------------------------------------
# 8. Leave the cursor where the data-entry form would naturally land
#    for the next new row.
# ------------------------------------------------------------------
[void]$ws.Range("P11").Select()
